$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in F1, copying the style used by the other header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 10:52:02.742781"
$ws.Range("F3").Value = "2021-10-05 10:52:02.742792"
$ws.Range("F4").Value = "2021-10-05 10:52:02.742796"
$ws.Range("F5").Value = "2021-10-05 10:52:02.742799"
$ws.Range("F6").Value = "2021-10-05 10:52:02.742801"
$ws.Range("F7").Value = "2021-10-05 10:52:02.742804"
$ws.Range("F8").Value = "2021-10-05 10:52:02.742807"
$ws.Range("F9").Value = "2021-10-05 10:52:02.742810"
$ws.Range("F10").Value = "2021-10-05 10:52:02.742813"
$ws.Range("F11").Value = "2021-10-05 10:52:02.742815"
$ws.Range("F12").Value = "2021-10-05 10:52:02.742818"
$ws.Range("F13").Value = "2021-10-05 10:52:02.742821"
